$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 61; existing rows 61-113 shift down to 62-114.
$ws.Rows.Item(61).Insert()

# Populate the new row 61 with the new weekly record.
$ws.Range("A61").Value = 7
$ws.Range("B61").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C61").Value = "Ñuble"
$ws.Range("D61").Value = 44942
$ws.Range("E61").Value = 16
$ws.Range("F61").Value = 100112031
$ws.Range("G61").Value = "Poroto verde"
$ws.Range("H61").Value = "Sin especificar"
$ws.Range("I61").Value = "Primera"
$ws.Range("J61").Value = 30
$ws.Range("K61").Value = 30000
$ws.Range("L61").Value = 30000
$ws.Range("M61").Value = 30000
$ws.Range("N61").Value = "$/saco 25 kilos"
$ws.Range("O61").Value = "Región del Maule"
$ws.Range("P61").Value = 1200
$ws.Range("Q61").Value = 25
$ws.Range("R61").Value = "Hortaliza"
